# Update correlation results for dataset 6 (DS 6) now that JaTyC
# (typestate_checker) warning counts have been fixed.
#
# Affects two sheets:
#   - "all_tools"         : rows 10-12 (dataset_id = 6)
#   - "typestate_checker"  : rows 10-12 (dataset_id = 6), plus a couple of
#                             column-width tweaks on columns I and K.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all_tools": F (num_snippets_warnings) 47 -> 48,
#                     G (num_warnings) 334 -> 819, and the recomputed
#                     correlation statistics in I/J/K/L for rows 10-12.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all_tools")

$allToolsRows = @{
    10 = @{ I = -0.01391413642584883; J = 0.9090032366864784; K = -0.01026612600739247; L = 0.9435901857829745 }
    11 = @{ I = -0.02815294687312959; J = 0.8086800956106934; K = -0.0299893827986723;  L = 0.8362130272367809 }
    12 = @{ I = -0.247884667760308;   J = 0.03188792894690582; K = -0.3096184720415943; L = 0.02866696147428332 }
}

foreach ($r in 10, 11, 12) {
    $wsAll.Cells.Item($r, 6).Value = 48    # F: num_snippets_warnings
    $wsAll.Cells.Item($r, 7).Value = 819   # G: num_warnings

    $vals = $allToolsRows[$r]
    $wsAll.Cells.Item($r, 9).Value = $vals.I   # I: kendalls_tau
    $wsAll.Cells.Item($r, 10).Value = $vals.J  # J: kendalls_p_value
    $wsAll.Cells.Item($r, 11).Value = $vals.K  # K: spearmans_rho
    $wsAll.Cells.Item($r, 12).Value = $vals.L  # L: spearmans_p_value
}

# ---------------------------------------------------------------------
# Sheet "typestate_checker": F (num_snippets_warnings) 5 -> 40,
#                             G (num_warnings) 52 -> 537, and the
#                             recomputed statistics for rows 10-12 (they
#                             now match the "all_tools" sheet since
#                             typestate_checker dominates the combined
#                             results for DS 6).
# ---------------------------------------------------------------------
$wsTypestate = $wb.Worksheets.Item("typestate_checker")

$typestateRows = @{
    10 = @{ I = -0.01391413642584883; J = 0.9090032366864784; K = -0.01026612600739247; L = 0.9435901857829745 }
    11 = @{ I = -0.02815294687312959; J = 0.8086800956106934; K = -0.0299893827986723;  L = 0.8362130272367809 }
    12 = @{ I = -0.247884667760308;   J = 0.03188792894690582; K = -0.3096184720415943; L = 0.02866696147428332 }
}

foreach ($r in 10, 11, 12) {
    $wsTypestate.Cells.Item($r, 6).Value = 40    # F: num_snippets_warnings
    $wsTypestate.Cells.Item($r, 7).Value = 537   # G: num_warnings

    $vals = $typestateRows[$r]
    $wsTypestate.Cells.Item($r, 9).Value = $vals.I   # I: kendalls_tau
    $wsTypestate.Cells.Item($r, 10).Value = $vals.J  # J: kendalls_p_value
    $wsTypestate.Cells.Item($r, 11).Value = $vals.K  # K: spearmans_rho
    $wsTypestate.Cells.Item($r, 12).Value = $vals.L  # L: spearmans_p_value
}

# Column width tweaks on "typestate_checker": columns I (9) and K (11)
# widen by one character (stored OOXML width 20.7109375 -> 21.7109375),
# matching the other sheets. The ColumnWidth COM property is in whole
# characters; Excel re-derives the fractional pixel-rounding remainder
# itself, so we set the integer character count (21) here -- the same
# width already used by column F (6) on this sheet and columns I/K on
# the other sheets.
$wsTypestate.Columns.Item(9).ColumnWidth = 21
$wsTypestate.Columns.Item(11).ColumnWidth = 21
